# Auto-generated edit script applying the Sheets diff (scheduled runner update).
# Updates cached leve-profit calculation values across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 563.86664
$ws.Range("J33").Value = 1233
$ws.Range("L33").Value = 1233
$ws.Range("N33").Value = -1691
$ws.Range("H41").Value = 8935.083000000001
$ws.Range("I41").Value = 1166
$ws.Range("J41").Value = 14484.429
$ws.Range("K41").Value = 1166
$ws.Range("L41").Value = 14484.429
$ws.Range("M41").Value = -726
$ws.Range("N41").Value = -15364.429
$ws.Range("H76").Value = 8649.799999999999
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 8649.799999999999
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H112").Value = 1973.75
$ws.Range("J112").Value = 1973.75
$ws.Range("L112").Value = 5921.25
$ws.Range("N112").Value = -8137.25
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27830.666
$ws.Range("I2").Value = 62546.2
$ws.Range("J2").Value = 3033.8572
$ws.Range("K2").Value = 62546.2
$ws.Range("L2").Value = 3033.8572
$ws.Range("M2").Value = -62433.2
$ws.Range("N2").Value = -3259.8572
$ws.Range("H88").Value = 2137
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2137
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2137
$ws.Range("N88").Value = -2949
$ws.Range("H91").Value = 2137
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2137
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2137
$ws.Range("N91").Value = -4945
$ws.Range("H116").Value = 27830.666
$ws.Range("I116").Value = 62546.2
$ws.Range("J116").Value = 3033.8572
$ws.Range("K116").Value = 62546.2
$ws.Range("L116").Value = 3033.8572
$ws.Range("M116").Value = -60252.2
$ws.Range("N116").Value = -7621.8572
$ws.Range("H122").Value = 4332.871
$ws.Range("I122").Value = 3890.4375
$ws.Range("J122").Value = 4804.8
$ws.Range("K122").Value = 11671.3125
$ws.Range("L122").Value = 14414.4
$ws.Range("M122").Value = -9221.3125
$ws.Range("N122").Value = -19314.4
$ws.Range("H132").Value = 2483.12
$ws.Range("I132").Value = 2169.9167
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 6509.750100000001
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3979.750100000001
$ws.Range("N132").Value = -35060
$ws.Range("M88").ClearContents()
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27830.666
$ws.Range("I3").Value = 62546.2
$ws.Range("J3").Value = 3033.8572
$ws.Range("K3").Value = 62546.2
$ws.Range("L3").Value = 3033.8572
$ws.Range("M3").Value = -62432.2
$ws.Range("N3").Value = -3261.8572
$ws.Range("H20").Value = 2154.05
$ws.Range("I20").Value = 1401.8334
$ws.Range("J20").Value = 3282.375
$ws.Range("K20").Value = 1401.8334
$ws.Range("L20").Value = 3282.375
$ws.Range("M20").Value = -1154.8334
$ws.Range("N20").Value = -3776.375
$ws.Range("H86").Value = 3900.1765
$ws.Range("I86").Value = 2678.5557
$ws.Range("J86").Value = 5274.5
$ws.Range("K86").Value = 2678.5557
$ws.Range("L86").Value = 5274.5
$ws.Range("M86").Value = -1555.5557
$ws.Range("N86").Value = -7520.5
$ws.Range("H89").Value = 3900.1765
$ws.Range("I89").Value = 2678.5557
$ws.Range("J89").Value = 5274.5
$ws.Range("K89").Value = 13392.7785
$ws.Range("L89").Value = 26372.5
$ws.Range("M89").Value = -7776.7785
$ws.Range("N89").Value = -37604.5
$ws.Range("H134").Value = 2190.8518
$ws.Range("I134").Value = 1806.12
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 5418.36
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -2883.36
$ws.Range("N134").Value = -26070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1105.7556
$ws.Range("I31").Value = 1144.4419
$ws.Range("K31").Value = 1144.4419
$ws.Range("M31").Value = -849.4419
$ws.Range("H34").Value = 1105.7556
$ws.Range("I34").Value = 1144.4419
$ws.Range("K34").Value = 1144.4419
$ws.Range("M34").Value = -942.4419
$ws.Range("H109").Value = 18659.076
$ws.Range("J109").Value = 18659.076
$ws.Range("L109").Value = 18659.076
$ws.Range("N109").Value = -20739.076
$ws.Range("H117").Value = 64500
$ws.Range("J117").Value = 64500
$ws.Range("L117").Value = 64500
$ws.Range("N117").Value = -73678
$ws.Range("H134").Value = 5390.6924
$ws.Range("I134").Value = 4505.4165
$ws.Range("K134").Value = 13516.2495
$ws.Range("M134").Value = -10981.2495

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1064.4445
$ws.Range("I5").Value = 963.75
$ws.Range("J5").Value = 1265.8334
$ws.Range("K5").Value = 2891.25
$ws.Range("L5").Value = 3797.5002
$ws.Range("M5").Value = -2779.25
$ws.Range("N5").Value = -4021.5002
$ws.Range("H12").Value = 2.909091
$ws.Range("I12").Value = 3
$ws.Range("K12").Value = 9
$ws.Range("M12").Value = 164
$ws.Range("H21").Value = 287.5
$ws.Range("I21").Value = 287.5
$ws.Range("K21").Value = 862.5
$ws.Range("M21").Value = -689.5
$ws.Range("H135").Value = 1064.4445
$ws.Range("I135").Value = 963.75
$ws.Range("J135").Value = 1265.8334
$ws.Range("K135").Value = 8673.75
$ws.Range("L135").Value = 11392.5006
$ws.Range("M135").Value = -6138.75
$ws.Range("N135").Value = -16462.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 10199.9
$ws.Range("I113").Value = 4200
$ws.Range("J113").Value = 16199.8
$ws.Range("K113").Value = 4200
$ws.Range("L113").Value = 16199.8
$ws.Range("M113").Value = -2030
$ws.Range("N113").Value = -20539.8
$ws.Range("H122").Value = 4854
$ws.Range("I122").Value = 3877.125
$ws.Range("K122").Value = 11631.375
$ws.Range("M122").Value = -9181.375
$ws.Range("H126").Value = 3213.125
$ws.Range("I126").Value = 2241.5715
$ws.Range("K126").Value = 6724.7145
$ws.Range("M126").Value = -4254.7145
$ws.Range("H132").Value = 6322.0713
$ws.Range("I132").Value = 6322.0713
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 18966.2139
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -16436.2139
$ws.Range("H134").Value = 80443.14
$ws.Range("J134").Value = 80443.14
$ws.Range("L134").Value = 241329.42
$ws.Range("N134").Value = -246399.42
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9560.866
$ws.Range("I40").Value = 11744.4
$ws.Range("J40").Value = 5193.8
$ws.Range("K40").Value = 11744.4
$ws.Range("L40").Value = 5193.8
$ws.Range("M40").Value = -11608.4
$ws.Range("N40").Value = -5465.8
$ws.Range("H64").Value = 60075
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 60075
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 60075
$ws.Range("N64").Value = -60525
$ws.Range("H67").Value = 60075
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 60075
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 60075
$ws.Range("N67").Value = -61635
$ws.Range("H136").Value = 3393.1304
$ws.Range("I136").Value = 3228.5789
$ws.Range("J136").Value = 4174.75
$ws.Range("K136").Value = 9685.736699999999
$ws.Range("L136").Value = 12524.25
$ws.Range("M136").Value = -7135.736699999999
$ws.Range("N136").Value = -17624.25
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13775.5
$ws.Range("I45").Value = 7968
$ws.Range("J45").Value = 15227.375
$ws.Range("K45").Value = 7968
$ws.Range("L45").Value = 15227.375
$ws.Range("M45").Value = -7477
$ws.Range("N45").Value = -16209.375
$ws.Range("H62").Value = 3568.5
$ws.Range("J62").Value = 3570
$ws.Range("L62").Value = 3570
$ws.Range("N62").Value = -4818
$ws.Range("H63").Value = 69249
$ws.Range("J63").Value = 69249
$ws.Range("L63").Value = 69249
$ws.Range("N63").Value = -70497
$ws.Range("H65").Value = 3568.5
$ws.Range("J65").Value = 3570
$ws.Range("L65").Value = 17850
$ws.Range("N65").Value = -24090
$ws.Range("H66").Value = 69249
$ws.Range("J66").Value = 69249
$ws.Range("L66").Value = 207747
$ws.Range("N66").Value = -213987
$ws.Range("H100").Value = 1522.8572
$ws.Range("I100").Value = 1709.8
$ws.Range("K100").Value = 3419.6
$ws.Range("M100").Value = -2878.6
$ws.Range("H101").Value = 21200.666
$ws.Range("J101").Value = 21200.666
$ws.Range("L101").Value = 21200.666
$ws.Range("N101").Value = -27690.666
$ws.Range("H107").Value = 991.5454999999999
$ws.Range("I107").Value = 942.6667
$ws.Range("K107").Value = 2828.0001
$ws.Range("M107").Value = -908.0001000000002
$ws.Range("H109").Value = 38318.184
$ws.Range("J109").Value = 38318.184
$ws.Range("L109").Value = 38318.184
$ws.Range("N109").Value = -41092.184
$ws.Range("H122").Value = 2247.8572
$ws.Range("I122").Value = 2052.2727
$ws.Range("J122").Value = 2965
$ws.Range("K122").Value = 6156.8181
$ws.Range("L122").Value = 8895
$ws.Range("M122").Value = -3706.8181
$ws.Range("N122").Value = -13795
$ws.Range("H136").Value = 1615.1875
$ws.Range("I136").Value = 1424.5714
$ws.Range("J136").Value = 2949.5
$ws.Range("K136").Value = 4273.7142
$ws.Range("L136").Value = 8848.5
$ws.Range("M136").Value = -1723.7142
$ws.Range("N136").Value = -13948.5
